$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency Price (D) / Volume(1h) (E) figures for this run.
# D-column "Price" values that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as literal text (preserving exact
# formatting such as trailing zeros, e.g. "1.00", "0.130", "2.60") instead of
# silently auto-converting them to numeric values.

$ws.Range("D2").Value = '66.516.05'
$ws.Range("E2").Value = '  +0.82%  '

$ws.Range("D3").Value = '3.280.17'
$ws.Range("E3").Value = '  +3.45%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = "'616.24"
$ws.Range("E5").Value = '  +1.98%  '

$ws.Range("D6").Value = "'159.05"
$ws.Range("E6").Value = '  +3.30%  '

$ws.Range("D8").Value = '3.280.57'
$ws.Range("E8").Value = '  +3.43%  '

$ws.Range("D9").Value = "'0.547"
$ws.Range("E9").Value = '  +0.60%  '

$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = '  +3.86%  '

$ws.Range("E11").Value = '  +2.26%  '

$ws.Range("D12").Value = "'0.499"
$ws.Range("E12").Value = '  -3.29%  '

$ws.Range("D13").Value = "'0.0000274"
$ws.Range("E13").Value = '  +3.69%  '

$ws.Range("D14").Value = "'39.26"
$ws.Range("E14").Value = '  +2.98%  '

$ws.Range("D15").Value = '3.806.61'
$ws.Range("E15").Value = '  +3.08%  '

$ws.Range("D16").Value = '66.591.21'
$ws.Range("E16").Value = '  +0.81%  '

$ws.Range("E17").Value = '  +1.29%  '

$ws.Range("D18").Value = '3.273.35'
$ws.Range("E18").Value = '  +3.08%  '

$ws.Range("E19").Value = '  +1.53%  '

$ws.Range("D20").Value = "'507.55"
$ws.Range("E20").Value = '  -0.03%  '

$ws.Range("E21").Value = '  +1.80%  '

$ws.Range("D22").Value = "'0.758"
$ws.Range("E22").Value = '  +4.38%  '

$ws.Range("E23").Value = '  +2.88%  '

$ws.Range("D24").Value = "'14.77"
$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("D25").Value = "'87.05"
$ws.Range("E25").Value = '  +3.27%  '

$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("E27").Value = '  +2.43%  '

$ws.Range("D28").Value = "'9.32"
$ws.Range("E28").Value = '  +2.24%  '

$ws.Range("E29").Value = '  +2.23%  '

$ws.Range("D30").Value = "'0.130"
$ws.Range("E30").Value = '  +48.70%  '

$ws.Range("D31").Value = "'7.04"

$ws.Range("D32").Value = "'2.90"
$ws.Range("E32").Value = '  -2.95%  '

$ws.Range("D33").Value = "'28.18"
$ws.Range("E33").Value = '  +1.08%  '

$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = '  -0.28%  '

$ws.Range("D35").Value = "'1.15"
$ws.Range("E35").Value = '  -1.91%  '

$ws.Range("D36").Value = "'6.52"
$ws.Range("E36").Value = '  +0.67%  '

$ws.Range("D37").Value = "'3.45"
$ws.Range("E37").Value = '  +20.97%  '

$ws.Range("E38").Value = '  +17.30%  '

$ws.Range("D39").Value = "'55.78"
$ws.Range("E39").Value = '  +0.84%  '

$ws.Range("D40").Value = "'496.77"
$ws.Range("E40").Value = '  -2.17%  '

$ws.Range("D41").Value = "'0.0427"
$ws.Range("E41").Value = '  +2.80%  '

$ws.Range("E42").Value = '  +1.81%  '

$ws.Range("E43").Value = '  +1.33%  '

$ws.Range("D44").Value = "'2.60"
$ws.Range("E44").Value = '  +6.84%  '

$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = "'0.296"
$ws.Range("E45").Value = '  +0.16%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '3.017.91'
$ws.Range("E46").Value = '  +6.87%  '

$ws.Range("D47").Value = "'29.26"
$ws.Range("E47").Value = '  +5.07%  '

$ws.Range("E48").Value = '  +6.48%  '

$ws.Range("E49").Value = '  +3.05%  '

$ws.Range("E50").Value = '  -0.02%  '

$ws.Range("D51").Value = "'2.55"
$ws.Range("E51").Value = '  -2.07%  '
